$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Sprint 2 "

$ws.Range("A2").Value = "1.-Antes de iniciar el modulo se debe configurar por estación , es decir que se debe identificar si es la primera vez que se abre  abrir un modal y configurar que numero de estación y guardarlo se me ocurre en el web.config"
$ws.Range("A3").Value = "2.-Se debe de poder guardar una compra, sin afectar inventario ni nada solo guardar una compra por si después desean venderla."
$ws.Range("A4").Value = "3.- de debe poder editar una compra aun que esta ya este finalizada (debemos pensar la forma de regresar producto al inventario o a ese venta agregar producto y marcar la salida de inventario pero solo afectar el inventario del producto que se esta agregando.)"
$ws.Range("A5").Value = "4.- En el modulo de ventas es importante tener en cuenta  que dependiendo de la cantidad del producto se aplica un precio y si el tipo de cliente es diferente a cliente general se aplica un ajuste al total de la venta. Esta por definirse el criterio de las 12 piezas el cual consiste en asignar un precio si el cliente lleva 12 piezas aunque no sea del mismo producto "
$ws.Range("A6").Value = "5.- debe poder buscar un producto y conocer sus distintos tipo de precios"
$ws.Range("A7").Value = "6.- Adicionar un modulo para configurar los rangos de precios"
$ws.Range("A8").Value = "7.-Adicionar un modulo para cargar el porcentaje de descuento dependiendo el tipo de cliente"
$ws.Range("A9").Value = "8.- Se debe realizar la facturacion al finalizar una compra  o despues de realizarla cuando se ven en el datatables"
$ws.Range("A1").Value = "Actividades "
$ws.Range("A10").Value = "9.-Agregar en el apartado de productos  un combo para seleccionar claveProdServ con base al catalogo FactCatClaveProdServicio"
$ws.Range("A11").Value = "10.-agregar en el apartado de productos un combo para seleccionar claveUnidadSaT con base al catalofo FactCatClaveUnidad"
$ws.Range("A12").Value = "11.- Generar ticket con la venta pensar que es un ticker como el de cualquier tienda chico no tamaño carta o algo asi como el que hizo blanquita para remesas "
$ws.Range("A13").Value = "12.-Agregar modulo para agregar estaciones por punto de venta (Seleccionar Sucursal , despues punto de venta , y ahí agregar estacion)"

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("A2:A13").Select()
